# Apply edits described in the commit diff across three sheets:
#   foresatt (sheet1), barn (sheet3), soknad (sheet4)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "foresatt"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("foresatt")

# Row 2
$ws.Range("C2").Value = "knsjndkfnj"
$ws.Range("D2").Value = "nfknkfndk"
$ws.Range("E2").Value = 232121232
$ws.Range("F2").Value = 434345

# Row 3
$ws.Range("C3").Value = "snnsjfndjs"
$ws.Range("D3").Value = "nkfnkfndk"
$ws.Range("E3").Value = 122232
$ws.Range("F3").Value = 121213342

# Row 4
$ws.Range("C4").Value = "Schlorpt"
$ws.Range("E4").Value = 8374363
$ws.Range("F4").Value = 12345678901

# Row 5
$ws.Range("C5").Value = "Glorpo"
$ws.Range("D5").Value = "Tilted Towers"
$ws.Range("E5").Value = 75757
$ws.Range("F5").Value = 100922846373

# Row 6
$ws.Range("C6").Value = "Schlorpt"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = 757574
$ws.Range("F6").Value = 230323434543

# Row 7
$ws.Range("C7").Value = "Glorpo"
$ws.Range("D7").Value = "Tilted Towers"
$ws.Range("E7").Value = 9865463
$ws.Range("F7").Value = 100922846373

# ---------------------------------------------------------------
# Sheet "barn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("barn")

# Row 2
$ws.Range("C2").Value = 343434

# Row 3
$ws.Range("C3").Value = 10070467433

# Row 4
$ws.Range("C4").Value = 10070467433

# ---------------------------------------------------------------
# Sheet "soknad"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("soknad")

# Row 2
$ws.Range("C2").Value = 5
$ws.Range("M2").Value = 23232

# Row 3
$ws.Range("H3").Value = ""
$ws.Range("J3").Value = "2,5,3"
$ws.Range("M3").Value = 21232

# Row 4
$ws.Range("G4").Value = ""
$ws.Range("J4").Value = "1,3,4,5"
$ws.Range("K4").Value = ""
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = 736372
